$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.488.07"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3
$ws.Range("D3").Value = "1.876.91"
$ws.Range("E3").Value = "  +1.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7161"
$ws.Range("E5").Value = "  +2.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.80"
$ws.Range("E6").Value = "  +1.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07978"
$ws.Range("E8").Value = "  +1.20%  "

# Row 9
$ws.Range("E9").Value = "  +3.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.37"
$ws.Range("E10").Value = "  +6.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08284"
$ws.Range("E11").Value = "  +1.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7304"
$ws.Range("E12").Value = "  +3.54%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.283"
$ws.Range("E13").Value = "  +1.90%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.858.15"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  +2.00%  "

# Row 16
$ws.Range("D16").Value = "29.492.55"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.932"
$ws.Range("E17").Value = "  +2.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "245.50"
$ws.Range("E18").Value = "  +4.18%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007879"
$ws.Range("E19").Value = "  +0.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("E20").Value = "  +1.05%  "

# Row 21
$ws.Range("D21").Value = "2.117.79"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.978"
$ws.Range("E23").Value = "  +6.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1610"
$ws.Range("E25").Value = "  +13.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.53"
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.058"
$ws.Range("E27").Value = "  +2.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("E28").Value = "  +1.58%  "

# Row 29
$ws.Range("E29").Value = "  -2.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.498"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.389"
$ws.Range("E31").Value = "  +2.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.119"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05268"
$ws.Range("E33").Value = "  +2.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.960"
$ws.Range("E34").Value = "  +2.28%  "

# Row 35
$ws.Range("E35").Value = "  +2.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7278"
$ws.Range("E36").Value = "  +3.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +1.39%  "

# Row 39
$ws.Range("D39").Value = "1.224.24"
$ws.Range("E39").Value = "  +6.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.715"
$ws.Range("E40").Value = "  +0.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9123"
$ws.Range("E41").Value = "  -1.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.87"
$ws.Range("E42").Value = "  +5.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.125"
$ws.Range("E43").Value = "  +3.00%  "

# Row 44
$ws.Range("E44").Value = "  +0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.25"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5286"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").Value = "2.014.14"
$ws.Range("E47").Value = "  +1.00%  "

# Row 48
$ws.Range("E48").Value = "  +3.94%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.931"
$ws.Range("E49").Value = "  +10.02%  "

# Row 50
$ws.Range("E50").Value = "  +1.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.339"
$ws.Range("E51").Value = "  +2.03%  "
